# feat: add 2022-Q1 data
#
# The workbook currently has two sheets: "2021-Q2" and "总计" (summary).
# We insert a brand-new "2022-Q1" sheet (holdings detail, same layout as
# "2021-Q2" but with an extra column) positioned between them, and add a
# corresponding "2022-Q1" row to the "总计" summary sheet (as the new first
# data row, pushing the existing "2021-Q2" row down).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by duplicating "总计" (so it starts
#    out with the same header/data-row cell styling), then placing the
#    duplicate right after "2021-Q2" (i.e. before "总计").
# ------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$summarySheet.Copy($summarySheet, $null) # insert copy immediately Before 总计

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"

# Extend the header formatting (style copied from 总计's D1) across the
# additional columns E1:H1 this sheet needs.
$newSheet.Range("D1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row 2 - B2:G2 are text cells (even though numeric-looking), H2 is
# a genuine number, matching the source data's typing.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "377016"
$newSheet.Range("C2").Value = "上投摩根亚太优势混合(QDII)"
$newSheet.Range("D2").Value = "27.15"
$newSheet.Range("E2").Value = "86.02"
$newSheet.Range("F2").Value = "1.83"
$newSheet.Range("G2").Value = "0.4968"
$newSheet.Range("B2:G2").ClearFormats()
$newSheet.Range("H2").Value = 8

# ------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing data row down to row 3
#    and write the new "2022-Q1" totals into row 2.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")   # same sheet object, now the third tab

# Carry A2's style down to A3 before either cell's value changes.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.78

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.5

# Restore the original active sheet/selection ("2021-Q2" was active before
# this edit and the sheet-copy operation above shifts focus to the new tab).
$wb.Worksheets.Item("2021-Q2").Activate()
